$d = $word.ActiveDocument

# 1. Remove the whole paragraph that contains the "{% load docx_tags %}"
#    template tag (it was split across three runs in the original markup).
#    Walk backwards so deleting doesn't shift the indices we still need to
#    examine.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*{% load*") {
        $p.Range.Delete()
    }
}

# 2. Flip "Allow punctuation to extend past text extents" (w:overflowPunct)
#    off for the Normal style - it was "true", it should become "false".
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false

Write-Output "done"
